{"js": "// Replace the 100 arithmetic-problem cell values in the single table\n// with the new values, preserving run/paragraph formatting (fonts,\n// size, alignment) by replacing text via the paragraph's Range rather\n// than the cell body (which would reset formatting to defaults).\n//\n// The table is 20 rows x 5 columns; the new value for each (row, col)\n// is taken from `newValues` below (row-major order, matching the\n// document/diff order).\n\nconst newValues = [\n  [\"99-12=\", \"17-10=\", \"68-16=\", \"35+7=\", \"27+31=\"],\n  [\"66-41=\", \"94-60=\", \"61-4=\", \"64+20=\", \"45+16=\"],\n  [\"8-4=\", \"40+5=\", \"68-56=\", \"66-11=\", \"50-2=\"],\n  [\"9+75=\", \"73+5=\", \"14+45=\", \"68-46=\", \"82-30=\"],\n  [\"34-15=\", \"79-31=\", \"9+21=\", \"6+41=\", \"49-22=\"],\n  [\"56+5=\", \"55-34=\", \"86-39=\", \"41+7=\", \"62-30=\"],\n  [\"77-24=\", \"7+81=\", \"44+14=\", \"45-43=\", \"92+2=\"],\n  [\"90-24=\", \"39+37=\", \"43-22=\", \"1+20=\", \"79-36=\"],\n  [\"65+20=\", \"28+53=\", \"36+6=\", \"79-56=\", \"32+3=\"],\n  [\"15+14=\", \"36-15=\", \"21+57=\", \"21-20=\", \"21-6=\"],\n  [\"57+8=\", \"13+48=\", \"97-77=\", \"59+8=\", \"92-32=\"],\n  [\"27+11=\", \"18-1=\", \"43+3=\", \"97-0=\", \"82-61=\"],\n  [\"59-36=\", \"29+14=\", \"12-5=\", \"31+35=\", \"69+15=\"],\n  [\"49-14=\", \"13+6=\", \"35+41=\", \"31+20=\", \"13+41=\"],\n  [\"23-0=\", \"33+16=\", \"65-59=\", \"12+22=\", \"65+6=\"],\n  [\"20-13=\", \"90+7=\", \"61+5=\", \"0+15=\", \"1+51=\"],\n  [\"99-50=\", \"20+48=\", \"97-97=\", \"40+21=\", \"37+50=\"],\n  [\"40+25=\", \"41+58=\", \"61-6=\", \"92-12=\", \"70+19=\"],\n  [\"40+35=\", \"77-38=\", \"33+38=\", \"64-35=\", \"57+29=\"],\n  [\"84-81=\", \"99-4=\", \"35+25=\", \"36-33=\", \"44+53=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rowCount = newValues.length;\nconst colCount = newValues[0].length;\n\n// Load all cell bodies' paragraphs in one go.\nconst cellParas = [];\nfor (let r = 0; r < rowCount; r++) {\n  const rowParas = [];\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const paras = cell.body.paragraphs;\n    paras.load(\"items\");\n    rowParas.push(paras);\n  }\n  cellParas.push(rowParas);\n}\nawait context.sync();\n\n// Replace the text of the (first/only) paragraph in each cell using\n// its Range, which preserves existing run formatting (font, size) and\n// paragraph formatting (alignment) instead of resetting to defaults.\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const para = cellParas[r][c].items[0];\n    const range = para.getRange();\n    range.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem cell values in the single table\n# with the new values, preserving run/paragraph formatting (fonts,\n# size, alignment). Setting Cell.Range.Text keeps the existing run\n# properties of the cell's paragraph mark, so formatting survives.\n#\n# The table is 20 rows x 5 columns; newValues[r][c] holds the\n# replacement text for that cell (row-major order, matching the\n# document/diff order).\n\n$newValues = @(\n    @(\"99-12=\", \"17-10=\", \"68-16=\", \"35+7=\", \"27+31=\"),\n    @(\"66-41=\", \"94-60=\", \"61-4=\", \"64+20=\", \"45+16=\"),\n    @(\"8-4=\", \"40+5=\", \"68-56=\", \"66-11=\", \"50-2=\"),\n    @(\"9+75=\", \"73+5=\", \"14+45=\", \"68-46=\", \"82-30=\"),\n    @(\"34-15=\", \"79-31=\", \"9+21=\", \"6+41=\", \"49-22=\"),\n    @(\"56+5=\", \"55-34=\", \"86-39=\", \"41+7=\", \"62-30=\"),\n    @(\"77-24=\", \"7+81=\", \"44+14=\", \"45-43=\", \"92+2=\"),\n    @(\"90-24=\", \"39+37=\", \"43-22=\", \"1+20=\", \"79-36=\"),\n    @(\"65+20=\", \"28+53=\", \"36+6=\", \"79-56=\", \"32+3=\"),\n    @(\"15+14=\", \"36-15=\", \"21+57=\", \"21-20=\", \"21-6=\"),\n    @(\"57+8=\", \"13+48=\", \"97-77=\", \"59+8=\", \"92-32=\"),\n    @(\"27+11=\", \"18-1=\", \"43+3=\", \"97-0=\", \"82-61=\"),\n    @(\"59-36=\", \"29+14=\", \"12-5=\", \"31+35=\", \"69+15=\"),\n    @(\"49-14=\", \"13+6=\", \"35+41=\", \"31+20=\", \"13+41=\"),\n    @(\"23-0=\", \"33+16=\", \"65-59=\", \"12+22=\", \"65+6=\"),\n    @(\"20-13=\", \"90+7=\", \"61+5=\", \"0+15=\", \"1+51=\"),\n    @(\"99-50=\", \"20+48=\", \"97-97=\", \"40+21=\", \"37+50=\"),\n    @(\"40+25=\", \"41+58=\", \"61-6=\", \"92-12=\", \"70+19=\"),\n    @(\"40+35=\", \"77-38=\", \"33+38=\", \"64-35=\", \"57+29=\"),\n    @(\"84-81=\", \"99-4=\", \"35+25=\", \"36-33=\", \"44+53=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
